# Apply updated coin data (price/volume/name/link refresh) per commit
# "Updated symbol list on Wed Jan 25 19:38:50 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text (these columns are stored as plain text,
# e.g. '301.41' not the number 301.41, matching the sheet's inlineStr cells).
$q = "'"

$ws.Range("D2").Value = $q + "301.41"
$ws.Range("E2").Value = $q + "-2.95%"
$ws.Range("D3").Value = $q + "35.51"
$ws.Range("E3").Value = $q + "-0.31%"
$ws.Range("D4").Value = $q + "5.072"
$ws.Range("E4").Value = $q + "-0.72%"
$ws.Range("D5").Value = $q + "0.08006"
$ws.Range("E5").Value = $q + "-2.63%"
$ws.Range("E6").Value = $q + "-6.37%"
$ws.Range("D7").Value = $q + "7.771"
$ws.Range("E7").Value = $q + "-2.11%"
$ws.Range("B8").Value = $q + "GateToken"
$ws.Range("C8").Value = $q + "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = $q + "4.042"
$ws.Range("E8").Value = $q + "-2.10%"
$ws.Range("B9").Value = $q + "MXToken"
$ws.Range("C9").Value = $q + "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = $q + "0.9282"
$ws.Range("E9").Value = $q + "0.13%"
$ws.Range("B10").Value = $q + "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = $q + "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = $q + "0.1549"
$ws.Range("E10").Value = $q + "37.12%"
$ws.Range("B11").Value = $q + "WazirX"
$ws.Range("C11").Value = $q + "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = $q + "0.1900"
$ws.Range("E11").Value = $q + "-0.89%"
$ws.Range("B12").Value = $q + "MandalaExchangeToken"
$ws.Range("C12").Value = $q + "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = $q + "0.09019"
$ws.Range("E12").Value = $q + "-3.29%"
$ws.Range("B13").Value = $q + "BitrueCoin"
$ws.Range("C13").Value = $q + "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = $q + "0.03457"
$ws.Range("E13").Value = $q + "-4.62%"
$ws.Range("B14").Value = $q + "BitMartToken"
$ws.Range("C14").Value = $q + "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = $q + "0.09885"
$ws.Range("E14").Value = $q + "-0.31%"
$ws.Range("B15").Value = $q + "BitForexToken"
$ws.Range("C15").Value = $q + "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = $q + "0.001399"
$ws.Range("E15").Value = $q + "-3.23%"
$ws.Range("B16").Value = $q + "TigerCash"
$ws.Range("C16").Value = $q + "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = $q + "0.005735"
$ws.Range("E16").Value = $q + "-1.60%"
$ws.Range("B17").Value = $q + "LEO"
$ws.Range("C17").Value = $q + "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = $q + "3.536"
$ws.Range("E17").Value = $q + "1.97%"
$ws.Range("E18").Value = $q + "0.09%"
$ws.Range("E19").Value = $q + "0.50%"
$ws.Range("D20").Value = $q + "0.1303"
$ws.Range("E20").Value = $q + "-0.58%"
$ws.Range("D21").Value = $q + "5.034"
$ws.Range("E21").Value = $q + "-1.18%"
$ws.Range("D23").Value = $q + "0.04486"
$ws.Range("E23").Value = $q + "-1.46%"
$ws.Range("D24").Value = $q + "0.001214"
$ws.Range("E24").Value = $q + "-1.01%"
$ws.Range("D25").Value = $q + "0.004773"
$ws.Range("E25").Value = $q + "-0.83%"
$ws.Range("D26").Value = $q + "0.0001231"
$ws.Range("E26").Value = $q + "-1.48%"
$ws.Range("D27").Value = $q + "0.0003023"
$ws.Range("D39").Value = $q + "0.01848"
$ws.Range("E39").Value = $q + "-6.67%"
$ws.Range("D40").Value = $q + "0.04776"
$ws.Range("E40").Value = $q + "-3.00%"
$ws.Range("E41").Value = $q + "6.02%"
$ws.Range("D42").Value = $q + "0.007335"
$ws.Range("E42").Value = $q + "-3.81%"
$ws.Range("D43").Value = $q + "0.1331"
$ws.Range("E43").Value = $q + "-3.92%"
$ws.Range("D44").Value = $q + "0.002112"
$ws.Range("E44").Value = $q + "-0.77%"
$ws.Range("D45").Value = $q + "0.009711"
$ws.Range("E45").Value = $q + "-16.05%"
$ws.Range("D46").Value = $q + "0.00006239"
$ws.Range("E46").Value = $q + "-4.73%"
$ws.Range("D47").Value = $q + "0.00000000751"
$ws.Range("E47").Value = $q + "0.11%"
$ws.Range("E48").Value = $q + "-63.72%"
$ws.Range("E49").Value = $q + "10.72%"
$ws.Range("D50").Value = $q + "0.00002102"
$ws.Range("E50").Value = $q + "0.11%"
$ws.Range("D51").Value = $q + "0.0002002"
$ws.Range("E51").Value = $q + "0.11%"
